$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NaukriSearch")

# Update the selection/active cell (cosmetic, matches diff)
$ws.Range("I9").Select()

# A2 = 1
$ws.Range("A2").Value = 1

# H2: "1 Month" -> "15 days" (new shared string, must be added before D2's edit below
# so it lands at shared-string index 34, matching the target ordering)
$ws.Range("H2").Value = "15 days"

# D2: "Asp.net, C#, Angular,SQL,Angular,MVC,Java" -> "Asp.net, C#, Angular,SQL,MVC,Java"
$ws.Range("D2").Value = "Asp.net, C#, Angular,SQL,MVC,Java"

# I2: 20 -> 2
$ws.Range("I2").Value = 2
